$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11 = 17.0306
    12 = 18.36510000000001
    15 = 16.4543
    27 = 16.76849999999998
    28 = 16.6292
    31 = 16.4873
    32 = 16.74829999999999
    36 = 17.29860000000001
    38 = 16.50429999999999
    46 = 17.16859999999999
    54 = 16.6966
    55 = 16.5556
    56 = 16.65360000000001
    67 = 17.15560000000003
    69 = 17.19190000000003
    72 = 16.46129999999999
    73 = 16.78840000000001
    83 = 16.77009999999999
    86 = 16.60750000000001
    91 = 18.72800000000002
    93 = 17.27290000000001
    99 = 16.3891
}

foreach ($row in $updates.Keys) {
    $ws.Range("E$row").Value = $updates[$row]
}
